$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.423.48'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.865.67'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -1.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.41'
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("E6").Value = '  -1.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5063'
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3888'
$ws.Range("E8").Value = '  -1.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08312'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +0.96%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.169'
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").Value = '1.864.38'
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.26'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.223'
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.97'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  -1.35%  '
$ws.Range("E22").Value = '  -1.08%  '
$ws.Range("D23").Value = '28.459.88'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.193'
$ws.Range("E25").Value = '  -4.41%  '
$ws.Range("D26").Value = '2.077.97'
$ws.Range("E26").Value = '  +1.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.77'
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.47'
$ws.Range("E28").Value = '  -1.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.398'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.43'
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1033'
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.032'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.755'
$ws.Range("E33").Value = '  -0.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.626'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06543'
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.929'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2152'
$ws.Range("E38").Value = '  -1.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.006'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6328'
$ws.Range("E42").Value = '  -1.22%  '
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5957'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.01'
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.674'
$ws.Range("E47").Value = '  -1.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.989'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.84'
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.205'
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.145'
$ws.Range("E51").Value = '  -6.26%  '
